$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B109").Value = "결제 타입"
$ws.Range("H109").Value = "결제 타입(1 : 무통장 입금, 2 : 실시간 계좌이체, 3 : 토스 페이)"
$ws.Range("B110").Value = "결제 금액"
$ws.Range("H110").Value = "결제 금액"
$ws.Range("B111").Value = "결제 금액 통화 코드"
$ws.Range("H111").Value = "결제 금액 통화 코드(IOS 4217, ex : KRW, USD, EUR...)"
$ws.Range("B112").Value = "결제 완료 일시"
$ws.Range("H112").Value = "결제가 완료 및 확인 된 일시(Null 이라면 아직 완료 처리가 아님)"
$ws.Range("B113").Value = "결제 실패 여부"
$ws.Range("H113").Value = "결제4 실패 여부"
$ws.Range("B120").Value = "결제 정보 고유키"
$ws.Range("B130").Value = "결제 정보 고유키"
$ws.Range("B140").Value = "결제 정보 고유키"
$ws.Range("B158").Value = "대여 가능 상품 예약 결제 정보 고유키"
$ws.Range("B159").Value = "결제 정보 고유키"
